# Generate Report for Handoff
# Renames the two tracked files, moves their status from
# "Handed back: in sync with en-US" to "Ready for handoff", regenerates the
# zh-cn/de-de xliff handoff rows (new handoff file + datetime, handback
# info cleared), and marks the second file as a content duplicate of the
# first (its handoff metadata mirrors row 2's newly generated xliff).

$wb = $excel.ActiveWorkbook

$oldName1 = "6a81eca0-262a-4276-91c6-6afdc564ae7b.md"
$newName1 = "57f738f8-fada-4700-8649-83148c0867a0.md"
$oldName2 = "b8051b9e-404c-4aa7-acd6-ac3064f8eac5.md"
$newName2 = "ffff2090f858-e132-4bdd-ae11-b9874b967283.md"

$newStatus = "Ready for handoff"
$newHoDate = "2016-08-27 23:01:36"

$newZhXlf = "57f738f8-fada-4700-8649-83148c0867a0.79a5a110f43a5efdcd6b1813773f0934f9958326.zh-cn.xlf"
$newDeXlf = "57f738f8-fada-4700-8649-83148c0867a0.79a5a110f43a5efdcd6b1813773f0934f9958326.de-de.xlf"
$newZhHandoffDate = "2016-08-27 23:01:32"
$clearedHandback = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newName1
$wsOverview.Range("A3").Value = $newName2

$wsOverview.Range("B2").Value = "e2e\" + $newName1
$wsOverview.Range("B3").Value = "e2e\" + $newName2

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsOverview.Range("G2").Value = $newHoDate
$wsOverview.Range("G3").Value = $newHoDate

foreach ($hl in $wsOverview.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq $wsOverview.Range("B2").Address()) {
        $hl.TextToDisplay = "e2e\" + $newName1
    } elseif ($addr -eq $wsOverview.Range("B3").Address()) {
        $hl.TextToDisplay = "e2e\" + $newName2
    }
}

$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333332
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333332

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newName1
$wsZh.Range("A3").Value = $newName2

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsZh.Range("F3").Value = "True"

$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("G3").Value = $newZhXlf

$wsZh.Range("H2").Value = $newZhHandoffDate
$wsZh.Range("H3").Value = $newZhHandoffDate

$wsZh.Range("I2").Value = ""
$wsZh.Range("I3").Value = ""
$wsZh.Range("J2").Value = ""
$wsZh.Range("J3").Value = ""

$wsZh.Range("K2").Value = $clearedHandback
$wsZh.Range("K3").Value = $clearedHandback

$zhTargets = @($wsZh.Range("I2").Address(), $wsZh.Range("I3").Address())
$zhHlList = @($wsZh.Hyperlinks)
foreach ($hl in $zhHlList) {
    $addr = $hl.Range.Address()
    if ($zhTargets -contains $addr) {
        $hl.Delete()
    } elseif ($addr -eq $wsZh.Range("A2").Address()) {
        $hl.TextToDisplay = $newName1
    } elseif ($addr -eq $wsZh.Range("A3").Address()) {
        $hl.TextToDisplay = $newName2
    }
}

$wsZh.Columns.Item(3).ColumnWidth = 16.333333333333332
$wsZh.Columns.Item(9).ColumnWidth = 17.833333333333332
$wsZh.Columns.Item(10).ColumnWidth = 20.833333333333332

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newName1
$wsDe.Range("A3").Value = $newName2

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$wsDe.Range("F3").Value = "True"

$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("G3").Value = $newDeXlf

$wsDe.Range("H2").Value = $newHoDate
$wsDe.Range("H3").Value = $newHoDate

$wsDe.Range("I2").Value = ""
$wsDe.Range("I3").Value = ""
$wsDe.Range("J2").Value = ""
$wsDe.Range("J3").Value = ""

$wsDe.Range("K2").Value = $clearedHandback
$wsDe.Range("K3").Value = $clearedHandback

$deTargets = @($wsDe.Range("I2").Address(), $wsDe.Range("I3").Address())
$deHlList = @($wsDe.Hyperlinks)
foreach ($hl in $deHlList) {
    $addr = $hl.Range.Address()
    if ($deTargets -contains $addr) {
        $hl.Delete()
    } elseif ($addr -eq $wsDe.Range("A2").Address()) {
        $hl.TextToDisplay = $newName1
    } elseif ($addr -eq $wsDe.Range("A3").Address()) {
        $hl.TextToDisplay = $newName2
    }
}

$wsDe.Columns.Item(3).ColumnWidth = 16.333333333333332
$wsDe.Columns.Item(9).ColumnWidth = 17.833333333333332
$wsDe.Columns.Item(10).ColumnWidth = 20.833333333333332
